$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows for the "abg_hypercap_threshold" group (rows 2-12)
$ws.Range("C2").Value = 2384
$ws.Range("D2").Value = 10747
$ws.Range("E2").Value = 22.18
$ws.Range("F2").Value = 21.4
$ws.Range("G2").Value = 22.97

$ws.Range("C3").Value = 1595
$ws.Range("D3").Value = 10747
$ws.Range("E3").Value = 14.84
$ws.Range("F3").Value = 14.17
$ws.Range("G3").Value = 15.51

$ws.Range("C4").Value = 1259
$ws.Range("D4").Value = 10747
$ws.Range("E4").Value = 11.71
$ws.Range("F4").Value = 11.11
$ws.Range("G4").Value = 12.32

$ws.Range("D5").Value = 10747

$ws.Range("C6").Value = 982
$ws.Range("D6").Value = 10747
$ws.Range("E6").Value = 9.140000000000001
$ws.Range("F6").Value = 8.59
$ws.Range("G6").Value = 9.68

$ws.Range("C7").Value = 950
$ws.Range("D7").Value = 10747
$ws.Range("E7").Value = 8.84
$ws.Range("F7").Value = 8.300000000000001
$ws.Range("G7").Value = 9.380000000000001

$ws.Range("C8").Value = 942
$ws.Range("D8").Value = 10747
$ws.Range("E8").Value = 8.77
$ws.Range("F8").Value = 8.23
$ws.Range("G8").Value = 9.300000000000001

$ws.Range("C9").Value = 627
$ws.Range("D9").Value = 10747
$ws.Range("E9").Value = 5.83
$ws.Range("F9").Value = 5.39
$ws.Range("G9").Value = 6.28

$ws.Range("C10").Value = 520
$ws.Range("D10").Value = 10747
$ws.Range("E10").Value = 4.84
$ws.Range("F10").Value = 4.43

$ws.Range("C11").Value = 285
$ws.Range("D11").Value = 10747
$ws.Range("E11").Value = 2.65
$ws.Range("F11").Value = 2.35
$ws.Range("G11").Value = 2.96

$ws.Range("C12").Value = 172
$ws.Range("D12").Value = 10747
$ws.Range("E12").Value = 1.6
$ws.Range("F12").Value = 1.36
$ws.Range("G12").Value = 1.84

# Update rows for the "vbg_hypercap_threshold" group (rows 35-45)
$ws.Range("C35").Value = 4457
$ws.Range("D35").Value = 17542
$ws.Range("E35").Value = 25.41
$ws.Range("F35").Value = 24.76
$ws.Range("G35").Value = 26.05

$ws.Range("C36").Value = 2612
$ws.Range("D36").Value = 17542
$ws.Range("E36").Value = 14.89
$ws.Range("F36").Value = 14.36
$ws.Range("G36").Value = 15.42

$ws.Range("C37").Value = 2459
$ws.Range("D37").Value = 17542
$ws.Range("E37").Value = 14.02
$ws.Range("F37").Value = 13.5
$ws.Range("G37").Value = 14.53

$ws.Range("C38").Value = 1689
$ws.Range("D38").Value = 17542
$ws.Range("E38").Value = 9.630000000000001
$ws.Range("F38").Value = 9.19
$ws.Range("G38").Value = 10.06

$ws.Range("C39").Value = 1124
$ws.Range("D39").Value = 17542
$ws.Range("E39").Value = 6.41
$ws.Range("F39").Value = 6.05
$ws.Range("G39").Value = 6.77

# Row 40 now holds the "Symptom – General" category (previously row 41's label)
$ws.Range("B40").Value = "Symptom – General"
$ws.Range("C40").Value = 1106
$ws.Range("D40").Value = 17542
$ws.Range("E40").Value = 6.3
$ws.Range("F40").Value = 5.95
$ws.Range("G40").Value = 6.66

# Row 41 now holds the "Injuries & adverse effects" category (previously row 40's label)
$ws.Range("B41").Value = "Injuries & adverse effects"
$ws.Range("C41").Value = 1100
$ws.Range("D41").Value = 17542
$ws.Range("E41").Value = 6.27
$ws.Range("F41").Value = 5.91

$ws.Range("C42").Value = 1065
$ws.Range("D42").Value = 17542
$ws.Range("F42").Value = 5.72

$ws.Range("C43").Value = 994
$ws.Range("D43").Value = 17542
$ws.Range("E43").Value = 5.67
$ws.Range("F43").Value = 5.32
$ws.Range("G43").Value = 6.01

$ws.Range("C44").Value = 556
$ws.Range("D44").Value = 17542
$ws.Range("E44").Value = 3.17
$ws.Range("F44").Value = 2.91
$ws.Range("G44").Value = 3.43

$ws.Range("C45").Value = 380
$ws.Range("D45").Value = 17542
$ws.Range("E45").Value = 2.17
$ws.Range("F45").Value = 1.95
$ws.Range("G45").Value = 2.38
